# Expand abbreviated gender codes in column B ("M" -> "Male", "F" -> "Female")
# across the data rows of the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -eq "M") {
        $cell.Value2 = "Male"
    } elseif ($val -eq "F") {
        $cell.Value2 = "Female"
    }
}
